$d = $word.ActiveDocument

# Turn on Track Changes so the three appended pieces of text land in their
# own runs (matching identically-formatted text normally gets merged into
# the neighboring run). We accept each tracked insertion individually
# afterwards so the final document has plain runs (no <w:ins> markup) while
# keeping the three new runs distinct from each other and from the
# pre-existing run.
$d.TrackRevisions = $true

# Locate the end of the first paragraph's text ("This is a Microsoft word document.")
$para = $d.Paragraphs(1)
$r = $para.Range
# Exclude the trailing paragraph mark so new text is appended to the same paragraph
$r.End = $r.End - 1
$r.Collapse(0)

$r.InsertAfter(" (")
$r.Collapse(0)

$r.InsertAfter("Changed main")
$r.Collapse(0)

$r.InsertAfter(")")

$d.TrackRevisions = $false

# Accept each tracked insertion one at a time (rather than
# AcceptAllRevisions, which forces a full-document repagination and wipes
# unrelated <w:lastRenderedPageBreak/> render hints elsewhere in the file).
for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
    $d.Revisions.Item($i).Accept()
}
